{"js": "// CIV-17609: updated GA documents to display main claim number.\n// The template previously labelled the case-number placeholders as\n// \"Claim number\"; both occurrences are relabelled to \"Case number\"\n// (the run formatting / surrounding merge fields are left untouched).\nconst body = context.document.body;\n\nconst results = body.search(\"Claim number\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Case number\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# CIV-17609: updated GA documents to display main claim number.\n# The template previously labelled the case-number placeholders as\n# \"Claim number\"; both occurrences are relabelled to \"Case number\"\n# (run formatting / surrounding merge fields are left untouched).\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Claim number\"\n$find.Replacement.Text = \"Case number\"\n\n$find.Execute(\n    $find.Text,              # FindText\n    $true,                   # MatchCase\n    $true,                   # MatchWholeWord\n    $false,                  # MatchWildcards\n    $false,                  # MatchSoundsLike\n    $false,                  # MatchAllWordForms\n    $true,                   # Forward\n    $wdFindContinue,         # Wrap\n    $false,                  # Format\n    $find.Replacement.Text,  # ReplaceWith\n    $wdReplaceAll            # Replace\n) | Out-Null\n"}
